$d = $word.ActiveDocument

# The date line reads "Kuni: 03.10.2023zh". The edit simply corrects the
# day digit from "3" to "4" (03.10.2023 -> 04.10.2023). A real editor did
# this by selecting the single "3" character and retyping "4", which is
# why the surrounding text ends up split into three runs with identical
# formatting ("Kuni: 0" | "4" | ".10.2023zh") instead of being re-joined
# into one run. We reproduce that precisely: briefly enable Track Changes
# for the one-character retype (so the engine keeps the run boundaries
# instead of silently re-merging same-formatted text on save), then accept
# the resulting revision right away so no tracked-change markup remains.

$d.TrackRevisions = $true

$rng = $d.Content
$found = $rng.Find.Execute("Күні: 03", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0) | Out-Null
    $rng.MoveStart(1, -1) | Out-Null   # extend selection back over the "3"
    $rng.Text = "4"

    $d.TrackRevisions = $false

    $para = $rng.Paragraphs(1).Range
    $revs = $para.Revisions
    for ($i = $revs.Count; $i -ge 1; $i--) {
        $revs.Item($i).Accept()
    }
} else {
    $d.TrackRevisions = $false
}
